$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Octubre de 2020 a las 08:38"

# Row 28 - Ucrania
$ws.Range("B28").Value = 276177
$ws.Range("C28").Value = 5590
$ws.Range("D28").Value = 118699
$ws.Range("E28").Value = 152249
$ws.Range("G28").Value = 107
$ws.Range("H28").Value = 5229

# Row 60 - Uzbekistan
$ws.Range("B60").Value = 61859
$ws.Range("C60").Value = 217
$ws.Range("E60").Value = 2584
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 512

# Row 68 - Kirguistan
$ws.Range("B68").Value = 50201
$ws.Range("C68").Value = 330
$ws.Range("D68").Value = 44884
$ws.Range("E68").Value = 4223
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 1094

# Row 75 - Afganistan
$ws.Range("B75").Value = 39994
$ws.Range("C75").Value = 66
$ws.Range("D75").Value = 33354
$ws.Range("E75").Value = 5160
